$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MotionPlanner")
$ws.Rows.Item(12).Delete()
$ws.Range("C3").Value = 10
